# Update "想去人数" (want-to-go count) values on the "展览" and "全部类型" sheets
# to reflect newly generated output (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 5998
$wsExhibit.Range("F6").Value = 112
$wsExhibit.Range("F10").Value = 31

# Sheet "全部类型" (all types) - mirrors the same events
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 5998
$wsAll.Range("F7").Value = 112
$wsAll.Range("F12").Value = 31
